# Updated symbol list on Tue Jan  3 19:48:45 UTC 2023 with GitHub Actions
# Refresh the crypto "Price" (D) and "Volume(1h)" (E) columns with the
# latest scraped values. The source values are textual (e.g. "245.21",
# "-0.46%") rather than numeric, so force each touched cell to Text
# format before writing the new value - this prevents Excel from
# re-interpreting them as numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

$updates = @(
    @{ Row = 2;  D = "245.21";     E = "-0.46%" },
    @{ Row = 3;  D = "28.44";      E = "-4.62%" },
    @{ Row = 4;  D = "5.236";      E = "1.21%" },
    @{ Row = 5;  D = "0.05700";    E = "-0.52%" },
    @{ Row = 6;  D = "6.612";      E = "0.53%" },
    @{ Row = 7;  E = "3.16%" },
    @{ Row = 8;  D = "0.8503";     E = "-0.68%" },
    @{ Row = 9;  D = "0.8671";     E = "-0.08%" },
    @{ Row = 10; D = "0.1368";     E = "0.15%" },
    @{ Row = 11; D = "0.07033";    E = "-0.47%" },
    @{ Row = 12; D = "0.03143";    E = "7.27%" },
    @{ Row = 13; D = "0.09204";    E = "-1.90%" },
    @{ Row = 14; D = "0.001524";   E = "0.67%" },
    @{ Row = 15; D = "0.0005979";  E = "-94.19%" },
    @{ Row = 16; D = "0.005941";   E = "-3.23%" },
    @{ Row = 17; D = "3.494";      E = "0.15%" },
    @{ Row = 18; E = "-4.42%" },
    @{ Row = 19; E = "0.45%" },
    @{ Row = 20; D = "0.03288";    E = "-4.57%" },
    @{ Row = 21; D = "0.1298";     E = "-0.97%" },
    @{ Row = 22; D = "3.516";      E = "1.25%" },
    @{ Row = 23; E = "-1.85%" },
    @{ Row = 25; E = "-0.33%" },
    @{ Row = 26; D = "0.004141";   E = "-17.43%" },
    @{ Row = 27; D = "0.0001199";  E = "-0.84%" },
    @{ Row = 40; D = "0.03767";    E = "0.36%" },
    @{ Row = 41; E = "-0.66%" },
    @{ Row = 42; D = "0.003735";   E = "-35.14%" },
    @{ Row = 43; D = "0.002199";   E = "4.74%" },
    @{ Row = 44; D = "0.009188";   E = "-4.17%" },
    @{ Row = 45; D = "0.00005274"; E = "0.57%" },
    @{ Row = 46; E = "-0.02%" },
    @{ Row = 48; D = "0.002438";   E = "-3.32%" },
    @{ Row = 49; D = "0.00002099"; E = "-0.02%" },
    @{ Row = 50; D = "0.0001999";  E = "-0.02%" }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey("D")) {
        Set-TextValue "D$row" $u.D
    }
    if ($u.ContainsKey("E")) {
        Set-TextValue "E$row" $u.E
    }
}
